# Update Tnfsf13-Sdc2 LR-pair TPM-derived statistics (columns E-T, rows 2-10)
# with new values computed from the updated TPM recomputation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6848073333333332
$ws.Range("H2").Value = 2.054422
$ws.Range("I2").Value = 0.2268310526442471
$ws.Range("J2").Value = 0.2268310526442472
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 0.3563378980162222
$ws.Range("R2").Value = 3.207041082146
$ws.Range("S2").Value = 0.0009313534166567637
$ws.Range("T2").Value = 0.0009313534166567641
$ws.Range("G3").Value = 0.6848073333333332
$ws.Range("H3").Value = 2.054422
$ws.Range("I3").Value = 0.2268310526442471
$ws.Range("J3").Value = 0.2268310526442472
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 72.56439445672643
$ws.Range("R3").Value = 653.0795501105379
$ws.Range("S3").Value = 0.189660143030379
$ws.Range("T3").Value = 0.189660143030379
$ws.Range("G4").Value = 0.6848073333333332
$ws.Range("H4").Value = 2.054422
$ws.Range("I4").Value = 0.2268310526442471
$ws.Range("J4").Value = 0.2268310526442472
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 13.86533516644
$ws.Range("R4").Value = 124.78801649796
$ws.Range("S4").Value = 0.0362395561972114
$ws.Range("T4").Value = 0.03623955619721141
$ws.Range("I5").Value = 0.1086184939966157
$ws.Range("J5").Value = 0.1086184939966157
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 0.1706331006502222
$ws.Range("R5").Value = 1.535697905852
$ws.Range("S5").Value = 0.0004459804083990167
$ws.Range("T5").Value = 0.0004459804083990168
$ws.Range("I6").Value = 0.1086184939966157
$ws.Range("J6").Value = 0.1086184939966157
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("S6").Value = 0.09081913109776753
$ws.Range("T6").Value = 0.09081913109776754
$ws.Range("I7").Value = 0.1086184939966157
$ws.Range("J7").Value = 0.1086184939966157
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 6.639442911280001
$ws.Range("R7").Value = 59.75498620152
$ws.Range("S7").Value = 0.01735338249044913
$ws.Range("T7").Value = 0.01735338249044913
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.006290666666667
$ws.Range("H8").Value = 6.018872
$ws.Range("I8").Value = 0.6645504533591371
$ws.Range("J8").Value = 0.6645504533591372
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 1.043968667055111
$ws.Range("R8").Value = 9.395718003496
$ws.Range("S8").Value = 0.002728600551210866
$ws.Range("T8").Value = 0.002728600551210867
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.006290666666667
$ws.Range("H9").Value = 6.018872
$ws.Range("I9").Value = 0.6645504533591371
$ws.Range("J9").Value = 0.6645504533591372
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 212.5930320024542
$ws.Range("R9").Value = 1913.337288022088
$ws.Range("S9").Value = 0.5556502628970792
$ws.Range("T9").Value = 0.5556502628970793
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.006290666666667
$ws.Range("H10").Value = 6.018872
$ws.Range("I10").Value = 0.6645504533591371
$ws.Range("J10").Value = 0.6645504533591372
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 40.62148750544
$ws.Range("R10").Value = 365.59338754896
$ws.Range("S10").Value = 0.106171589910847
$ws.Range("T10").Value = 0.1061715899108471
